$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$wsf = $excel.WorksheetFunction

# Rows 9, 10, 11, 13, 14, 15, 16 have their entire record (all populated
# columns A:AY) permuted between each other; row 12 is left untouched.
# Map: destination row -> source row (source values are the *original*,
# pre-edit contents, since several rows are both a source and a
# destination in this permutation, e.g. row 9's old data ends up in row 11
# while row 9 itself receives row 16's old data).
$mapping = @{
    9  = 16
    10 = 15
    11 = 9
    13 = 11
    14 = 10
    15 = 14
    16 = 13
}

$lastCol = 51   # column AY - the last used column on this sheet

function Get-CellState($row, $col) {
    $cell = $ws.Cells.Item($row, $col)
    $blank = $wsf.IsBlank($cell)
    if ($blank) {
        return @{ Present = $false; Value = $null }
    } else {
        return @{ Present = $true; Value = $cell.Value2 }
    }
}

# 1) Snapshot every distinct source row's cells (presence + value) before
#    any writes happen, since rows get read after other rows have already
#    been overwritten otherwise.
$snapshot = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowStates = @()
        for ($col = 1; $col -le $lastCol; $col++) {
            $rowStates += , (Get-CellState $srcRow $col)
        }
        $snapshot[$srcRow] = $rowStates
    }
}

# 2) Apply the snapshotted data to each destination row, only touching
#    cells whose target state actually differs from what's already there
#    (keeps unaffected/constant columns byte-identical).
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowStates = $snapshot[$srcRow]
    for ($col = 1; $col -le $lastCol; $col++) {
        $target = $rowStates[$col - 1]
        $current = Get-CellState $destRow $col
        $cell = $ws.Cells.Item($destRow, $col)

        if ($target.Present) {
            if ((-not $current.Present) -or ($current.Value -ne $target.Value)) {
                $cell.Value = $target.Value
            }
        } else {
            if ($current.Present) {
                $cell.ClearContents()
            }
        }
    }
}
